$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rightmost column of each quartile group (K, then F) so the
# "Developed" block shrinks from B:F to B:E and the "Emerging" block shrinks
# from G:K to F:I.
$ws.Columns("K").Delete()
$ws.Columns("F").Delete()

# The eight quartile-average cells now live in B4:I4; overwrite with the
# newly-processed modeling data.
$ws.Cells.Item(4, 2).Value = 0.007686371481006832
$ws.Cells.Item(4, 3).Value = 0.007377612257013231
$ws.Cells.Item(4, 4).Value = 0.00828968790674171
$ws.Cells.Item(4, 5).Value = 0.00780329286272748
$ws.Cells.Item(4, 6).Value = 0.01407305369895101
$ws.Cells.Item(4, 7).Value = 0.01507321394877679
$ws.Cells.Item(4, 8).Value = 0.01087102909208167
$ws.Cells.Item(4, 9).Value = 0.006243591214958307
